# Insert a new row for the "LP solver (linprog or gurobi)" setting into the
# "general" sheet, right after the "NLP solver" row and before the
# "Number of exp. conditions" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("general")

# Push row 5 (and everything below) down by one row.
$ws.Rows.Item(5).Insert()

# Fill in the newly inserted row.
$ws.Cells.Item(5, 1).Value = "LP solver (linprog or gurobi)"
$ws.Cells.Item(5, 2).Value = "gurobi"
